# Update keyboard type names and refresh the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Find and replace the two shared-string values used in column C.
$found = $ws.Cells.Find("Dome-Switch")
if ($found) {
    $ws.Cells.Replace("Dome-Switch", "Dome", 1, 1, $false, $false, $false, $false)
}
$found2 = $ws.Cells.Find("Scissor-Switch")
if ($found2) {
    $ws.Cells.Replace("Scissor-Switch", "Scissor", 1, 1, $false, $false, $false, $false)
}

# Update the active selection shown in the sheet view to C2.
$ws.Activate()
$ws.Range("C2").Select()
